# Add newly-discovered identity elements info to the aaRS identity elements sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ArgRS: add a note in column C with additional identity element numbers
$ws.Range("C3").Value = "19, 56, 17"

# GlyRS: position 37 was added to the known identity elements, and
# highlighted separately in the new column C
$ws.Range("B9").Value = "1, 2, 3, 35, 36, 37, 70, 71, 72, 73"
$ws.Range("C9").Value = 37
